$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11 (Leve Item ID 5533)
$ws.Range("H11").Value = 67.14286
$ws.Range("I11").Value = 67.14286
$ws.Range("K11").Value = 67.14286
$ws.Range("M11").Value = 72.85714

# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 11539.2
$ws.Range("I86").Value = 1725
$ws.Range("K86").Value = 1725
$ws.Range("M86").Value = -602

# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 11539.2
$ws.Range("I89").Value = 1725
$ws.Range("K89").Value = 8625
$ws.Range("M89").Value = -3009

# Row 103 (Leve Item ID 19909)
$ws.Range("H103").Value = 66666908
$ws.Range("I103").Value = 111111250
$ws.Range("K103").Value = 333333750
$ws.Range("M103").Value = -333333164

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 2388.275
$ws.Range("I132").Value = 2411.0513
$ws.Range("K132").Value = 7233.1539
$ws.Range("M132").Value = -4703.1539

$ws = $wb.Worksheets.Item("ARM")
# Row 22 (Leve Item ID 2497)
$ws.Range("H22").Value = 1332.625
$ws.Range("I22").Value = 1332.625
$ws.Range("K22").Value = 1332.625
$ws.Range("M22").Value = -1033.625

# Row 44 (Leve Item ID 3861)
$ws.Range("H44").Value = 29088.75
$ws.Range("J44").Value = 29088.75
$ws.Range("L44").Value = 29088.75
$ws.Range("N44").Value = -30064.75

# Row 63 (Leve Item ID 12528)
$ws.Range("H63").Value = 2085557.9
$ws.Range("I63").Value = 2428.3845
$ws.Range("K63").Value = 2428.3845
$ws.Range("M63").Value = -1742.3845

# Row 66 (Leve Item ID 12528)
$ws.Range("H66").Value = 2085557.9
$ws.Range("I66").Value = 2428.3845
$ws.Range("K66").Value = 12141.9225
$ws.Range("M66").Value = -8709.922500000001

# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 545.1818
$ws.Range("I97").Value = 562.8946999999999
$ws.Range("J97").Value = 433
$ws.Range("K97").Value = 562.8946999999999
$ws.Range("L97").Value = 433
$ws.Range("M97").Value = -66.89469999999994
$ws.Range("N97").Value = -1425

# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 1745.6923
$ws.Range("I102").Value = 1724.5
$ws.Range("K102").Value = 1724.5
$ws.Range("M102").Value = -102.5

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 2180.8845
$ws.Range("I122").Value = 2105.348
$ws.Range("J122").Value = 2760
$ws.Range("K122").Value = 6316.044
$ws.Range("L122").Value = 8280
$ws.Range("M122").Value = -3866.044
$ws.Range("N122").Value = -13180

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 11851.607
$ws.Range("I132").Value = 1848.2285
$ws.Range("J132").Value = 33734
$ws.Range("K132").Value = 5544.6855
$ws.Range("L132").Value = 101202
$ws.Range("M132").Value = -3014.6855
$ws.Range("N132").Value = -106262

$ws = $wb.Worksheets.Item("BSM")
# Row 82 (Leve Item ID 11877)
$ws.Range("H82").Value = 31441.5

# Row 85 (Leve Item ID 11877)
$ws.Range("H85").Value = 31441.5

$ws = $wb.Worksheets.Item("CRP")
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 18522104
$ws.Range("I99").Value = 3075.9167
$ws.Range("K99").Value = 3075.9167
$ws.Range("M99").Value = -1577.9167

# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 855.5625
$ws.Range("I122").Value = 855.5625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2566.6875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -116.6875
$ws.Range("N122").ClearContents()

# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 18522104
$ws.Range("I126").Value = 3075.9167
$ws.Range("K126").Value = 9227.750100000001
$ws.Range("M126").Value = -6757.750100000001

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2466.4062
$ws.Range("I132").Value = 1830.5454
$ws.Range("J132").Value = 3865.3
$ws.Range("K132").Value = 5491.6362
$ws.Range("L132").Value = 11595.9
$ws.Range("M132").Value = -2961.6362
$ws.Range("N132").Value = -16655.9

# Row 137 (Leve Item ID 43231)
$ws.Range("H137").Value = 50780
$ws.Range("J137").Value = 50780
$ws.Range("L137").Value = 50780
$ws.Range("N137").Value = -60980

$ws = $wb.Worksheets.Item("CUL")
# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 880.6
$ws.Range("J122").Value = 976.9048
$ws.Range("L122").Value = 8792.1432
$ws.Range("N122").Value = -13692.1432

# Row 123 (Leve Item ID 36037)
$ws.Range("H123").Value = 4443
$ws.Range("I123").Value = 3499.25
$ws.Range("J123").Value = 5198
$ws.Range("K123").Value = 10497.75
$ws.Range("L123").Value = 15594
$ws.Range("M123").Value = -8047.75
$ws.Range("N123").Value = -20494

# Row 125 (Leve Item ID 36043)
$ws.Range("H125").Value = 4999.75
$ws.Range("J125").Value = 4999.75
$ws.Range("L125").Value = 14999.25
$ws.Range("N125").Value = -24839.25

# Row 129 (Leve Item ID 36054)
$ws.Range("H129").Value = 201788.9
$ws.Range("I129").Value = 582.5
$ws.Range("J129").Value = 252090.5
$ws.Range("K129").Value = 1747.5
$ws.Range("L129").Value = 756271.5
$ws.Range("M129").Value = 3252.5
$ws.Range("N129").Value = -766271.5

# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 644.33
$ws.Range("J131").Value = 768.0137
$ws.Range("L131").Value = 2304.0411
$ws.Range("N131").Value = -12384.0411

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Range("H2").Value = 85.46154
$ws.Range("I2").Value = 87.09999999999999
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 87.09999999999999
$ws.Range("L2").Value = 80
$ws.Range("M2").Value = 25.90000000000001
$ws.Range("N2").Value = -306

# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 2849479.5
$ws.Range("I70").Value = 4444.4443
$ws.Range("J70").Value = 4819119
$ws.Range("K70").Value = 4444.4443
$ws.Range("L70").Value = 4819119
$ws.Range("M70").Value = -4174.4443
$ws.Range("N70").Value = -4819659

# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 2849479.5
$ws.Range("I73").Value = 4444.4443
$ws.Range("J73").Value = 4819119
$ws.Range("K73").Value = 4444.4443
$ws.Range("L73").Value = 4819119
$ws.Range("M73").Value = -3508.4443
$ws.Range("N73").Value = -4820991

# Row 97 (Leve Item ID 19940)
$ws.Range("H97").Value = 2819.111
$ws.Range("I97").Value = 2669.6
$ws.Range("J97").Value = 3566.6667
$ws.Range("K97").Value = 2669.6
$ws.Range("L97").Value = 3566.6667
$ws.Range("M97").Value = -2173.6
$ws.Range("N97").Value = -4558.6667

# Row 114 (Leve Item ID 25957)
$ws.Range("H114").Value = 30000
$ws.Range("J114").Value = 30000
$ws.Range("L114").Value = 30000
$ws.Range("N114").Value = -38678

# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 6584.9287
$ws.Range("I122").Value = 7086.25
$ws.Range("K122").Value = 21258.75
$ws.Range("M122").Value = -18808.75

$ws = $wb.Worksheets.Item("LTW")
# Row 2 (Leve Item ID 2631)
$ws.Range("H2").Value = 498750
$ws.Range("I2").Value = 500000
$ws.Range("J2").Value = 490000
$ws.Range("K2").Value = 500000
$ws.Range("L2").Value = 490000
$ws.Range("M2").Value = -499888
$ws.Range("N2").Value = -490224

# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 2687.3125
$ws.Range("I7").Value = 2438.2307
$ws.Range("J7").Value = 3766.6667
$ws.Range("K7").Value = 2438.2307
$ws.Range("L7").Value = 3766.6667
$ws.Range("M7").Value = -2326.2307
$ws.Range("N7").Value = -3990.6667

# Row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 1401.6111
$ws.Range("I82").Value = 1389.3125
$ws.Range("J82").Value = 1500
$ws.Range("K82").Value = 1389.3125
$ws.Range("L82").Value = 1500
$ws.Range("M82").Value = -1028.3125
$ws.Range("N82").Value = -2222

# Row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 1401.6111
$ws.Range("I85").Value = 1389.3125
$ws.Range("J85").Value = 1500
$ws.Range("K85").Value = 1389.3125
$ws.Range("L85").Value = 1500
$ws.Range("M85").Value = -141.3125
$ws.Range("N85").Value = -3996

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 1229676.8
$ws.Range("I122").Value = 1638019
$ws.Range("J122").Value = 4650
$ws.Range("K122").Value = 4914057
$ws.Range("L122").Value = 13950
$ws.Range("M122").Value = -4911607
$ws.Range("N122").Value = -18850

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 2687.3125
$ws.Range("I126").Value = 2438.2307
$ws.Range("J126").Value = 3766.6667
$ws.Range("K126").Value = 7314.6921
$ws.Range("L126").Value = 11300.0001
$ws.Range("M126").Value = -4844.6921
$ws.Range("N126").Value = -16240.0001

# Row 141 (Leve Item ID 42487)
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

$ws = $wb.Worksheets.Item("WVR")
# Row 13 (Leve Item ID 3008)
$ws.Range("H13").Value = 2000
$ws.Range("J13").Value = 2000
$ws.Range("L13").Value = 2000
$ws.Range("N13").Value = -2280

# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 1772.5
$ws.Range("I81").Value = 696.6667
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 1393.3334
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -332.3334
$ws.Range("N81").Value = -12122

# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 1772.5
$ws.Range("I84").Value = 696.6667
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 6966.666999999999
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -1662.666999999999
$ws.Range("N84").Value = -60608

# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 87121730
$ws.Range("I107").Value = 111111550
$ws.Range("J107").Value = 15152255
$ws.Range("K107").Value = 333334650
$ws.Range("L107").Value = 45456765
$ws.Range("M107").Value = -333332730
$ws.Range("N107").Value = -45460605

# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 1288017.6
$ws.Range("J113").Value = 3003530
$ws.Range("L113").Value = 9010590
$ws.Range("N113").Value = -9014930

# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 1760
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050

# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 2140.8125
$ws.Range("I126").Value = 1562.3334
$ws.Range("J126").Value = 3876.25
$ws.Range("K126").Value = 4687.0002
$ws.Range("L126").Value = 11628.75
$ws.Range("M126").Value = -2217.0002
$ws.Range("N126").Value = -16568.75

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 620.3022999999999
$ws.Range("I132").Value = 489.16666
$ws.Range("K132").Value = 1467.49998
$ws.Range("M132").Value = 1062.50002
